$d = $word.ActiveDocument

# Remove the first paragraph entirely (including its paragraph mark) --
# the disclaimer text "I personally examined the patient separately ..."
$p = $d.Paragraphs(1)
$p.Range.Delete()
